$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append 4 new data rows (37-40) below the existing last row (36).
# Copy the formatting of the last existing data row for each new row so
# that number formats / styles (date, text, number) are preserved exactly
# as used throughout the rest of the table.

$ws.Rows("36").Copy()
$ws.Rows("37").Insert(-4121)

$ws.Rows("36").Copy()
$ws.Rows("38").Insert(-4121)

$ws.Rows("36").Copy()
$ws.Rows("39").Insert(-4121)

$ws.Rows("36").Copy()
$ws.Rows("40").Insert(-4121)

# Row 37: 04/06/2018 - Licata Rosa - Tela Leggera - Mt. - 62
$ws.Range("A37").Value = 43255
$ws.Range("B37").Value = "Licata Rosa"
$ws.Range("C37").Value = "Tela Leggera"
$ws.Range("D37").Value = "Mt."
$ws.Range("E37").Value = 62

# Row 38: 04/06/2018 - Licata Rosa - Silesia Bianca - Mt. - 60
$ws.Range("A38").Value = 43255
$ws.Range("B38").Value = "Licata Rosa"
$ws.Range("C38").Value = "Silesia Bianca"
$ws.Range("D38").Value = "Mt."
$ws.Range("E38").Value = 60

# Row 39: 04/06/2018 - Licata Rosa - Adesivo Strech Bianco - Mt. - 10
$ws.Range("A39").Value = 43255
$ws.Range("B39").Value = "Licata Rosa"
$ws.Range("C39").Value = "Adesivo Strech Bianco"
$ws.Range("D39").Value = "Mt."
$ws.Range("E39").Value = 10

# Row 40: 09/06/2018 - Bertolotti Daniela - Tela Leggera - Mt. - 7
$ws.Range("A40").Value = 43260
$ws.Range("B40").Value = "Bertolotti Daniela"
$ws.Range("C40").Value = "Tela Leggera"
$ws.Range("D40").Value = "Mt."
$ws.Range("E40").Value = 7
